$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '69.200.65'
$ws.Range('E2').Value = '  +0.16%  '
Set-TextValue 'D3' '3.754.42'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  -0.05%  '
Set-TextValue 'D5' '602.57'
$ws.Range('E5').Value = '  +0.07%  '
Set-TextValue 'D6' '167.47'
$ws.Range('E6').Value = '  -0.41%  '
Set-TextValue 'D7' '3.753.56'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +1.20%  '
Set-TextValue 'D10' '0.171'
$ws.Range('E10').Value = '  +2.71%  '
Set-TextValue 'D11' '6.40'
$ws.Range('E11').Value = '  +1.55%  '
$ws.Range('E12').Value = '  -0.01%  '
Set-TextValue 'D13' '38.03'
$ws.Range('E13').Value = '  -0.77%  '
Set-TextValue 'D15' '4.384.40'
$ws.Range('E15').Value = '  +0.13%  '
Set-TextValue 'D16' '3.750.78'
$ws.Range('E16').Value = '  +0.12%  '
Set-TextValue 'D17' '69.193.50'
$ws.Range('E17').Value = '  +0.20%  '
Set-TextValue 'D18' '7.36'
$ws.Range('E18').Value = '  +1.28%  '
Set-TextValue 'D19' '17.37'
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('E20').Value = '  -1.60%  '
Set-TextValue 'D21' '11.13'
$ws.Range('E21').Value = '  +6.85%  '
Set-TextValue 'D22' '493.74'
$ws.Range('E22').Value = '  -1.04%  '
Set-TextValue 'D23' '0.729'
$ws.Range('E23').Value = '  +0.51%  '
Set-TextValue 'D24' '0.0000153'
$ws.Range('E24').Value = '  +8.31%  '
Set-TextValue 'D25' '84.94'
$ws.Range('E25').Value = '  -0.26%  '
Set-TextValue 'D26' '2.30'
$ws.Range('E26').Value = '  -0.39%  '
Set-TextValue 'D27' '12.31'
$ws.Range('E27').Value = '  -0.04%  '
Set-TextValue 'D28' '10.10'
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('E30').Value = '  +1.07%  '
Set-TextValue 'D31' '8.14'
$ws.Range('E31').Value = '  +2.42%  '
Set-TextValue 'D32' '2.47'
$ws.Range('E32').Value = '  +2.24%  '
Set-TextValue 'D33' '31.58'
$ws.Range('E33').Value = '  -0.78%  '
Set-TextValue 'D34' '3.902.38'
$ws.Range('E34').Value = '  +0.29%  '
Set-TextValue 'D35' '3.688.37'
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('E36').Value = '  -0.60%  '
Set-TextValue 'D37' '1.00'
$ws.Range('E37').Value = '  -0.07%  '
Set-TextValue 'D38' '5.99'
$ws.Range('E38').Value = '  +3.06%  '
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('E40').Value = '  +3.57%  '
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('E42').Value = '  +5.69%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D43' '48.66'
$ws.Range('E43').Value = '  -0.89%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D44' '426.55'
$ws.Range('E44').Value = '  -3.33%  '
Set-TextValue 'D45' '1.98'
$ws.Range('E45').Value = '  -0.46%  '
$ws.Range('E46').Value = '  +0.74%  '
$ws.Range('E47').Value = '  -0.01%  '
Set-TextValue 'D48' '40.27'
$ws.Range('E48').Value = '  -0.40%  '
Set-TextValue 'D49' '141.57'
$ws.Range('E49').Value = '  -0.98%  '
Set-TextValue 'D50' '2.794.85'
$ws.Range('E50').Value = '  +1.64%  '
Set-TextValue 'D51' '0.0354'
$ws.Range('E51').Value = '  +0.63%  '
